$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H41").Value = 1721.6923
$ws.Range("I41").Value = 1667.1428
$ws.Range("J41").Value = 1785.3334
$ws.Range("K41").Value = 1667.1428
$ws.Range("L41").Value = 1785.3334
$ws.Range("M41").Value = -1227.1428
$ws.Range("N41").Value = -2665.3334

$ws.Range("H98").Value = 2790.1304
$ws.Range("I98").Value = 2428.2942
$ws.Range("J98").Value = 3815.3333
$ws.Range("K98").Value = 2428.2942
$ws.Range("L98").Value = 3815.3333
$ws.Range("M98").Value = -930.2941999999998

$ws.Range("H122").Value = 2790.1304
$ws.Range("I122").Value = 2428.2942
$ws.Range("J122").Value = 3815.3333
$ws.Range("K122").Value = 7284.882599999999
$ws.Range("L122").Value = 11445.9999
$ws.Range("M122").Value = -4834.882599999999

$ws.Range("H132").Value = 7577765.5
$ws.Range("I132").Value = 8132041.5
$ws.Range("J132").Value = 2660
$ws.Range("K132").Value = 24396124.5
$ws.Range("L132").Value = 7980
$ws.Range("M132").Value = -24393594.5
$ws.Range("N132").Value = -13040

$ws.Range("H138").Value = 1337.92
$ws.Range("I138").Value = 722.9268
$ws.Range("J138").Value = 1765.2881
$ws.Range("K138").Value = 2168.7804
$ws.Range("L138").Value = 5295.8643
$ws.Range("M138").Value = 2971.2196
$ws.Range("N138").Value = -15575.8643

$ws.Range("H140").Value = 30811.428
$ws.Range("I140").Value = 0
$ws.Range("J140").Value = 30811.428
$ws.Range("K140").Value = 0
$ws.Range("L140").Value = 30811.428
$ws.Range("N140").Value = -41171.428

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 15345.714
$ws.Range("I2").Value = 1176.75
$ws.Range("J2").Value = 34237.668
$ws.Range("K2").Value = 1176.75
$ws.Range("L2").Value = 34237.668
$ws.Range("M2").Value = -1063.75

$ws.Range("H61").Value = 2805.25
$ws.Range("I61").Value = 2222
$ws.Range("J61").Value = 2999.6667
$ws.Range("K61").Value = 2222
$ws.Range("L61").Value = 2999.6667
$ws.Range("M61").Value = -2010
$ws.Range("N61").Value = -3423.6667

$ws.Range("H97").Value = 323.75
$ws.Range("I97").Value = 323.75
$ws.Range("J97").Value = 0
$ws.Range("K97").Value = 323.75
$ws.Range("L97").Value = 0
$ws.Range("M97").Value = 172.25

$ws.Range("H116").Value = 15345.714
$ws.Range("I116").Value = 1176.75
$ws.Range("J116").Value = 34237.668
$ws.Range("K116").Value = 1176.75
$ws.Range("L116").Value = 34237.668
$ws.Range("M116").Value = 1117.25

$ws.Range("H136").Value = 2805.25
$ws.Range("I136").Value = 2222
$ws.Range("J136").Value = 2999.6667
$ws.Range("K136").Value = 6666
$ws.Range("L136").Value = 8999.000100000001
$ws.Range("M136").Value = -4116
$ws.Range("N136").Value = -14099.0001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 15345.714
$ws.Range("I3").Value = 1176.75
$ws.Range("J3").Value = 34237.668
$ws.Range("K3").Value = 1176.75
$ws.Range("L3").Value = 34237.668
$ws.Range("M3").Value = -1062.75

$ws.Range("H55").Value = 10709
$ws.Range("I55").Value = 10709
$ws.Range("J55").Value = 0
$ws.Range("K55").Value = 10709
$ws.Range("L55").Value = 0
$ws.Range("M55").Value = -10436

$ws.Range("H134").Value = 8910.883
$ws.Range("I134").Value = 953.2727
$ws.Range("J134").Value = 23499.834
$ws.Range("K134").Value = 2859.8181
$ws.Range("L134").Value = 70499.50199999999
$ws.Range("M134").Value = -324.8181
$ws.Range("N134").Value = -75569.50199999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 125001030
$ws.Range("I16").Value = 125001030
$ws.Range("J16").Value = 0
$ws.Range("K16").Value = 125001030
$ws.Range("L16").Value = 0
$ws.Range("M16").Value = -125000743

$ws.Range("H31").Value = 1407.1538
$ws.Range("I31").Value = 1117.6364
$ws.Range("J31").Value = 2999.5
$ws.Range("K31").Value = 1117.6364
$ws.Range("L31").Value = 2999.5
$ws.Range("M31").Value = -822.6364000000001

$ws.Range("H34").Value = 1407.1538
$ws.Range("I34").Value = 1117.6364
$ws.Range("J34").Value = 2999.5
$ws.Range("K34").Value = 1117.6364
$ws.Range("L34").Value = 2999.5
$ws.Range("M34").Value = -915.6364000000001

$ws.Range("H113").Value = 125001030
$ws.Range("I113").Value = 125001030
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 125001030
$ws.Range("L113").Value = 0
$ws.Range("M113").Value = -124998860

$ws.Range("H132").Value = 12787
$ws.Range("I132").Value = 18328.666
$ws.Range("J132").Value = 4474.5
$ws.Range("K132").Value = 54985.99800000001
$ws.Range("L132").Value = 13423.5
$ws.Range("M132").Value = -52455.99800000001

$ws.Range("H134").Value = 22224044
$ws.Range("I134").Value = 27779480
$ws.Range("J134").Value = 2300
$ws.Range("K134").Value = 83338440
$ws.Range("L134").Value = 6900
$ws.Range("M134").Value = -83335905

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 2339917.5
$ws.Range("I4").Value = 449744.5
$ws.Range("J4").Value = 3600033
$ws.Range("K4").Value = 1349233.5
$ws.Range("L4").Value = 10800099
$ws.Range("M4").Value = -1349121.5
$ws.Range("N4").Value = -10800323

$ws.Range("H34").Value = 1595.375
$ws.Range("I34").Value = 421.5
$ws.Range("J34").Value = 2299.7
$ws.Range("K34").Value = 1264.5
$ws.Range("L34").Value = 6899.099999999999
$ws.Range("M34").Value = -1180.5
$ws.Range("N34").Value = -7067.099999999999

$ws.Range("H47").Value = 0
$ws.Range("I47").Value = 0
$ws.Range("J47").Value = 0
$ws.Range("K47").Value = 0
$ws.Range("L47").Value = 0
$ws.Range("M47").ClearContents()

$ws.Range("H80").Value = 5071.4287
$ws.Range("I80").Value = 0
$ws.Range("J80").Value = 5071.4287
$ws.Range("K80").Value = 0
$ws.Range("L80").Value = 15214.2861
$ws.Range("N80").Value = -17086.2861

$ws.Range("H83").Value = 5071.4287
$ws.Range("I83").Value = 0
$ws.Range("J83").Value = 5071.4287
$ws.Range("K83").Value = 0
$ws.Range("L83").Value = 45642.85830000001
$ws.Range("N83").Value = -55002.85830000001

$ws.Range("H113").Value = 600
$ws.Range("I113").Value = 0
$ws.Range("J113").Value = 600
$ws.Range("K113").Value = 0
$ws.Range("L113").Value = 1800
$ws.Range("N113").Value = -6140
$ws.Range("M113").ClearContents()

$ws.Range("H121").Value = 426.4
$ws.Range("I121").Value = 359
$ws.Range("J121").Value = 1033
$ws.Range("K121").Value = 1077
$ws.Range("L121").Value = 3099
$ws.Range("M121").Value = 233
$ws.Range("N121").Value = -5719

$ws.Range("H132").Value = 1204
$ws.Range("I132").Value = 1204
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 10836
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -8306
$ws.Range("N132").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 3565.2727
$ws.Range("I80").Value = 2016
$ws.Range("J80").Value = 6276.5
$ws.Range("K80").Value = 2016
$ws.Range("L80").Value = 6276.5
$ws.Range("M80").Value = -1018
$ws.Range("N80").Value = -8272.5

$ws.Range("H82").Value = 19417.5
$ws.Range("I82").Value = 0
$ws.Range("J82").Value = 19417.5
$ws.Range("K82").Value = 0
$ws.Range("L82").Value = 19417.5
$ws.Range("N82").Value = -20183.5

$ws.Range("H83").Value = 3565.2727
$ws.Range("I83").Value = 2016
$ws.Range("J83").Value = 6276.5
$ws.Range("K83").Value = 10080
$ws.Range("L83").Value = 31382.5
$ws.Range("M83").Value = -5088
$ws.Range("N83").Value = -41366.5

$ws.Range("H85").Value = 19417.5
$ws.Range("I85").Value = 0
$ws.Range("J85").Value = 19417.5
$ws.Range("K85").Value = 0
$ws.Range("L85").Value = 19417.5
$ws.Range("N85").Value = -22069.5

$ws.Range("H126").Value = 2702.5
$ws.Range("I126").Value = 2505.8333
$ws.Range("J126").Value = 2997.5
$ws.Range("K126").Value = 7517.499899999999
$ws.Range("L126").Value = 8992.5
$ws.Range("M126").Value = -5047.499899999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H11").Value = 4000
$ws.Range("I11").Value = 0
$ws.Range("J11").Value = 4000
$ws.Range("K11").Value = 0
$ws.Range("L11").Value = 4000
$ws.Range("N11").Value = -4280

$ws.Range("H18").Value = 0
$ws.Range("I18").Value = 0
$ws.Range("J18").Value = 0
$ws.Range("K18").Value = 0
$ws.Range("L18").Value = 0
$ws.Range("M18").ClearContents()

$ws.Range("H61").Value = 1266.6666
$ws.Range("I61").Value = 1000
$ws.Range("J61").Value = 1400
$ws.Range("K61").Value = 1000
$ws.Range("L61").Value = 1400
$ws.Range("M61").Value = -798
$ws.Range("N61").Value = -1804

$ws.Range("H113").Value = 1266.6666
$ws.Range("I113").Value = 1000
$ws.Range("J113").Value = 1400
$ws.Range("K113").Value = 1000
$ws.Range("L113").Value = 1400
$ws.Range("M113").Value = 1170
$ws.Range("N113").Value = -5740

$ws.Range("H136").Value = 14622
$ws.Range("I136").Value = 51242
$ws.Range("J136").Value = 2415.3333
$ws.Range("K136").Value = 153726
$ws.Range("L136").Value = 7245.999899999999
$ws.Range("M136").Value = -151176
$ws.Range("N136").Value = -12345.9999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H12").Value = 70007
$ws.Range("I12").Value = 0
$ws.Range("J12").Value = 70007
$ws.Range("K12").Value = 0
$ws.Range("L12").Value = 70007
$ws.Range("N12").Value = -70291

$ws.Range("H17").Value = 1000
$ws.Range("I17").Value = 1000
$ws.Range("J17").Value = 0
$ws.Range("K17").Value = 1000
$ws.Range("L17").Value = 0
$ws.Range("M17").Value = -828
$ws.Range("N17").ClearContents()

$ws.Range("H81").Value = 850
$ws.Range("I81").Value = 850
$ws.Range("J81").Value = 0
$ws.Range("K81").Value = 1700
$ws.Range("L81").Value = 0
$ws.Range("M81").Value = -639

$ws.Range("H84").Value = 850
$ws.Range("I84").Value = 850
$ws.Range("J84").Value = 0
$ws.Range("K84").Value = 8500
$ws.Range("L84").Value = 0
$ws.Range("M84").Value = -3196

$ws.Range("H126").Value = 111112150
$ws.Range("I126").Value = 111112150
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 333336450
$ws.Range("L126").Value = 0
$ws.Range("M126").Value = -333333980
